{"js": "// Insert a new \"List Bullet\" paragraph with the docente's name right\n// after the \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"Docente(s) Respons\u00e1vel(eis)\"\n);\n\nif (!target) {\n  throw new Error('Could not find the \"Docente(s) Respons\u00e1vel(eis)\" paragraph');\n}\n\nconst newPara = target.insertParagraph(\"5701460 - Antonio Iacono\", \"After\");\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph with the docente's name right\n# after the \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw 'Could not find the \"Docente(s) Respons\u00e1vel(eis)\" paragraph'\n}\n\n$r = $target.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.InsertBefore(\"5701460 - Antonio Iacono\")\n$newPara.Style = \"List Bullet\"\n"}
